$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Word"
$ws.Range("B1").Value = "Meaning"

$ws.Range("A3").Value = "('Word', 'sound')"
$ws.Range("B3").Value = "('Meaning', 'goood')"
